$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1190.0435
$ws.Range("I15").Value = 1190.0435
$ws.Range("K15").Value = 3570.1305
$ws.Range("M15").Value = -3401.1305

# Row 17
$ws.Range("H17").Value = 610.0952
$ws.Range("J17").Value = 585.3889
$ws.Range("L17").Value = 1756.1667
$ws.Range("N17").Value = -2092.1667

# Row 33
$ws.Range("H33").Value = 7249.2666
$ws.Range("I33").Value = 43
$ws.Range("J33").Value = 9050.833000000001
$ws.Range("K33").Value = 43
$ws.Range("L33").Value = 9050.833000000001
$ws.Range("M33").Value = 186
$ws.Range("N33").Value = -9508.833000000001

# Row 40
$ws.Range("H40").Value = 1993.3334
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1990
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1990
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2340

# Row 64
$ws.Range("H64").Value = 3324.9592
$ws.Range("I64").Value = 3154.8572
$ws.Range("K64").Value = 3154.8572
$ws.Range("M64").Value = -2906.8572

# Row 67
$ws.Range("H67").Value = 3324.9592
$ws.Range("I67").Value = 3154.8572
$ws.Range("K67").Value = 3154.8572
$ws.Range("M67").Value = -2296.8572

# Row 106
$ws.Range("H106").Value = 2148.1765
$ws.Range("I106").Value = 1968.5
$ws.Range("K106").Value = 1968.5
$ws.Range("M106").Value = -1337.5

# Row 113
$ws.Range("H113").Value = 2284.353
$ws.Range("I113").Value = 2177.8333
$ws.Range("J113").Value = 2540
$ws.Range("K113").Value = 2177.8333
$ws.Range("L113").Value = 2540
$ws.Range("M113").Value = 1076.1667
$ws.Range("N113").Value = -9048

# Row 132
$ws.Range("H132").Value = 3332.9807
$ws.Range("I132").Value = 3620.75
$ws.Range("J132").Value = 1750.25
$ws.Range("K132").Value = 10862.25
$ws.Range("L132").Value = 5250.75
$ws.Range("M132").Value = -8332.25
$ws.Range("N132").Value = -10310.75

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 8333
$ws.Range("I63").Value = 7999
$ws.Range("J63").Value = 8399.799999999999
$ws.Range("K63").Value = 7999
$ws.Range("L63").Value = 8399.799999999999
$ws.Range("M63").Value = -7313
$ws.Range("N63").Value = -9771.799999999999

# Row 66
$ws.Range("H66").Value = 8333
$ws.Range("I66").Value = 7999
$ws.Range("J66").Value = 8399.799999999999
$ws.Range("K66").Value = 39995
$ws.Range("L66").Value = 41999
$ws.Range("M66").Value = -36563
$ws.Range("N66").Value = -48863

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1201.7778
$ws.Range("I99").Value = 754.73334
$ws.Range("J99").Value = 3437
$ws.Range("K99").Value = 754.73334
$ws.Range("L99").Value = 3437
$ws.Range("M99").Value = 743.26666
$ws.Range("N99").Value = -6433

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 285010.97
$ws.Range("I99").Value = 379425.8
$ws.Range("J99").Value = 1766.4667
$ws.Range("K99").Value = 379425.8
$ws.Range("L99").Value = 1766.4667
$ws.Range("M99").Value = -377927.8
$ws.Range("N99").Value = -4762.4667

# Row 126
$ws.Range("H126").Value = 285010.97
$ws.Range("I126").Value = 379425.8
$ws.Range("J126").Value = 1766.4667
$ws.Range("K126").Value = 1138277.4
$ws.Range("L126").Value = 5299.4001
$ws.Range("M126").Value = -1135807.4
$ws.Range("N126").Value = -10239.4001

# Row 132
$ws.Range("H132").Value = 36014.5
$ws.Range("I132").Value = 1736.8
$ws.Range("J132").Value = 104569.9
$ws.Range("K132").Value = 5210.4
$ws.Range("L132").Value = 313709.7
$ws.Range("M132").Value = -2680.4
$ws.Range("N132").Value = -318769.7

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 401.66666
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 401.66666
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2409.99996
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -2635.99996

# Row 38
$ws.Range("H38").Value = 1053.9032
$ws.Range("I38").Value = 662.4666999999999
$ws.Range("K38").Value = 1987.4001
$ws.Range("M38").Value = -1640.4001

# Row 122
$ws.Range("H122").Value = 464.2069
$ws.Range("I122").Value = 311.47827
$ws.Range("J122").Value = 1049.6666
$ws.Range("K122").Value = 2803.30443
$ws.Range("L122").Value = 9446.999400000001
$ws.Range("M122").Value = -353.3044300000001
$ws.Range("N122").Value = -14346.9994

# Row 123
$ws.Range("H123").Value = 1457.5
$ws.Range("I123").Value = 1457.5
$ws.Range("K123").Value = 4372.5
$ws.Range("M123").Value = -1922.5

# Row 131
$ws.Range("H131").Value = 12501218
$ws.Range("I131").Value = 3418
$ws.Range("J131").Value = 14286619
$ws.Range("K131").Value = 10254
$ws.Range("L131").Value = 42859857
$ws.Range("M131").Value = -5214
$ws.Range("N131").Value = -42869937

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 16551270
$ws.Range("I70").Value = 28131380
$ws.Range("J70").Value = 8257.714
$ws.Range("K70").Value = 28131380
$ws.Range("L70").Value = 8257.714
$ws.Range("M70").Value = -28131110
$ws.Range("N70").Value = -8797.714

# Row 73
$ws.Range("H73").Value = 16551270
$ws.Range("I73").Value = 28131380
$ws.Range("J73").Value = 8257.714
$ws.Range("K73").Value = 28131380
$ws.Range("L73").Value = 8257.714
$ws.Range("M73").Value = -28130444
$ws.Range("N73").Value = -10129.714

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# Row 87
$ws.Range("H87").Value = 22000
$ws.Range("I87").Value = 22000
$ws.Range("K87").Value = 22000
$ws.Range("M87").Value = -20877

# Row 88
$ws.Range("H88").Value = 2171
$ws.Range("I88").Value = 2171
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 2171
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1743
$ws.Range("N88").ClearContents()

# Row 90
$ws.Range("H90").Value = 22000
$ws.Range("I90").Value = 22000
$ws.Range("K90").Value = 66000
$ws.Range("M90").Value = -60384

# Row 91
$ws.Range("H91").Value = 2171
$ws.Range("I91").Value = 2171
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 2171
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -689
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 22938
$ws.Range("J75").Value = 22938
$ws.Range("L75").Value = 22938
$ws.Range("N75").Value = -24810

# Row 78
$ws.Range("H78").Value = 22938
$ws.Range("J78").Value = 22938
$ws.Range("L78").Value = 68814
$ws.Range("N78").Value = -78174

# Row 136
$ws.Range("H136").Value = 26139.05
$ws.Range("I136").Value = 35116.45
$ws.Range("J136").Value = 2471.3635
$ws.Range("K136").Value = 105349.35
$ws.Range("L136").Value = 7414.0905
$ws.Range("M136").Value = -102799.35
$ws.Range("N136").Value = -12514.0905
